$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-17 (columns A-T), reflecting the updated
# TPM-based NATMI ligand-receptor signalling table (Icosl -> Icos), now
# expanded to the full 4x4 sending/target cluster cross join.
$data = @(
    @("ECs", "Icosl", "Icos", "ECs", 3, 1, 0.9314806666666667, 2.794442, 0.06288211749152639, 0.06288211749152638, 1, 0.3333333333333333, 0.03464866666666667, 0.103946, 0.03773320526419342, 0.03773320526419342, 0.03227456312577778, 0.290471068132, 0.002372743846754893, 0.002372743846754892),
    @("ECs", "Icosl", "Icos", "FAPs", 3, 1, 0.9314806666666667, 2.794442, 0.06288211749152639, 0.06288211749152638, 2, 0.6666666666666666, 0.3099333333333333, 0.9298, 0.3375246210017417, 0.3375246210017416, 0.2886969079555556, 2.5982721716, 0.02122426287411444, 0.02122426287411443),
    @("ECs", "Icosl", "Icos", "MuSCs", 3, 1, 0.9314806666666667, 2.794442, 0.06288211749152639, 0.06288211749152638, 1, 0.3333333333333333, 0.02645, 0.07935, 0.02880466624702969, 0.02880466624702969, 0.02463766363333334, 0.2217389727, 0.001811298407249926, 0.001811298407249925),
    @("ECs", "Icosl", "Icos", "Resolving-Mac", 3, 1, 0.9314806666666667, 2.794442, 0.06288211749152639, 0.06288211749152638, 3, 1, 0.547222, 1.641666, 0.5959375074870351, 0.5959375074870351, 0.5097267133746667, 4.587540420372, 0.03747381236340713, 0.03747381236340712),
    @("FAPs", "Icosl", "Icos", "ECs", 3, 1, 2.399281333333333, 7.197844, 0.1619699646990985, 0.1619699646990985, 1, 0.3333333333333333, 0.03464866666666667, 0.103946, 0.03773320526419342, 0.03773320526419342, 0.08313189915822222, 0.7481870924239999, 0.006111645924625246, 0.006111645924625245),
    @("FAPs", "Icosl", "Icos", "FAPs", 3, 1, 2.399281333333333, 7.197844, 0.1619699646990985, 0.1619699646990985, 2, 0.6666666666666666, 0.3099333333333333, 0.9298, 0.3375246210017417, 0.3375246210017416, 0.7436172612444444, 6.692555351199999, 0.0546688509487287, 0.05466885094872869),
    @("FAPs", "Icosl", "Icos", "MuSCs", 3, 1, 2.399281333333333, 7.197844, 0.1619699646990985, 0.1619699646990985, 1, 0.3333333333333333, 0.02645, 0.07935, 0.02880466624702969, 0.02880466624702969, 0.06346099126666667, 0.5711489214000001, 0.004665490775200713, 0.004665490775200712),
    @("FAPs", "Icosl", "Icos", "Resolving-Mac", 3, 1, 2.399281333333333, 7.197844, 0.1619699646990985, 0.1619699646990985, 3, 1, 0.547222, 1.641666, 0.5959375074870351, 0.5959375074870351, 1.312939529789333, 11.816455768104, 0.09652397705054382, 0.09652397705054382),
    @("MuSCs", "Icosl", "Icos", "ECs", 3, 1, 0.416663, 1.249989, 0.02812796084553394, 0.02812796084553394, 1, 0.3333333333333333, 0.03464866666666667, 0.103946, 0.03773320526419342, 0.03773320526419342, 0.01443681739933333, 0.129931356594, 0.001061358120247728, 0.001061358120247727),
    @("MuSCs", "Icosl", "Icos", "FAPs", 3, 1, 0.416663, 1.249989, 0.02812796084553394, 0.02812796084553394, 2, 0.6666666666666666, 0.3099333333333333, 0.9298, 0.3375246210017417, 0.3375246210017416, 0.1291377524666666, 1.1622397722, 0.009493879323940673, 0.009493879323940671),
    @("MuSCs", "Icosl", "Icos", "MuSCs", 3, 1, 0.416663, 1.249989, 0.02812796084553394, 0.02812796084553394, 1, 0.3333333333333333, 0.02645, 0.07935, 0.02880466624702969, 0.02880466624702969, 0.01102073635, 0.09918662714999998, 0.0008102165243651242, 0.0008102165243651241),
    @("MuSCs", "Icosl", "Icos", "Resolving-Mac", 3, 1, 0.416663, 1.249989, 0.02812796084553394, 0.02812796084553394, 3, 1, 0.547222, 1.641666, 0.5959375074870351, 0.5959375074870351, 0.228007160186, 2.052064441674, 0.01676250687698042, 0.01676250687698041),
    @("Resolving-Mac", "Icosl", "Icos", "ECs", 3, 1, 11.0657, 33.1971, 0.7470199569638412, 0.7470199569638412, 1, 0.3333333333333333, 0.03464866666666667, 0.103946, 0.03773320526419342, 0.03773320526419342, 0.3834117507333333, 3.4507057566, 0.02818745737256555, 0.02818745737256555),
    @("Resolving-Mac", "Icosl", "Icos", "FAPs", 3, 1, 11.0657, 33.1971, 0.7470199569638412, 0.7470199569638412, 2, 0.6666666666666666, 0.3099333333333333, 0.9298, 0.3375246210017417, 0.3375246210017416, 3.429629286666667, 30.86666358, 0.2521376278549579, 0.2521376278549579),
    @("Resolving-Mac", "Icosl", "Icos", "MuSCs", 3, 1, 11.0657, 33.1971, 0.7470199569638412, 0.7470199569638412, 1, 0.3333333333333333, 0.02645, 0.07935, 0.02880466624702969, 0.02880466624702969, 0.292687765, 2.634189885, 0.02151766054021393, 0.02151766054021393),
    @("Resolving-Mac", "Icosl", "Icos", "Resolving-Mac", 3, 1, 11.0657, 33.1971, 0.7470199569638412, 0.7470199569638412, 3, 1, 0.547222, 1.641666, 0.5959375074870351, 0.5959375074870351, 6.0553944854, 54.4985503686, 0.4451772111961038, 0.4451772111961038)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($excelRow, $j + 1).Value = $rowValues[$j]
    }
}
